$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 4760
$ws.Range('L3').Value = 5137
$ws.Range('L4').Value = 1261
$ws.Range('L5').Value = 302
$ws.Range('L6').Value = 4357
$ws.Range('L7').Value = 15817

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 131
$ws.Range('L4').Value = 60
$ws.Range('L5').Value = 57
$ws.Range('L7').Value = 521
$ws.Range('L8').Value = 1053
$ws.Range('L11').Value = 258
$ws.Range('L14').Value = 84
$ws.Range('L19').Value = 434
$ws.Range('L20').Value = 396
$ws.Range('L21').Value = 50
$ws.Range('L23').Value = 173
$ws.Range('L24').Value = 42
$ws.Range('L29').Value = 866
$ws.Range('L30').Value = 76
$ws.Range('L33').Value = 726
$ws.Range('L37').Value = 583
$ws.Range('L41').Value = 72
$ws.Range('L48').Value = 204
$ws.Range('L50').Value = 77
$ws.Range('L51').Value = 195
$ws.Range('L52').Value = 317
$ws.Range('L53').Value = 180
$ws.Range('L54').Value = 331
$ws.Range('L55').Value = 151
$ws.Range('L65').Value = 306
$ws.Range('L67').Value = 543
$ws.Range('L73').Value = 124
$ws.Range('L75').Value = 57
$ws.Range('L76').Value = 245
$ws.Range('L79').Value = 417
$ws.Range('L83').Value = 351
$ws.Range('L85').Value = 813
$ws.Range('L86').Value = 115
$ws.Range('L89').Value = 227
$ws.Range('L91').Value = 216
$ws.Range('L93').Value = 83
$ws.Range('L96').Value = 175
$ws.Range('L101').Value = 15817

$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('L3').Value = 19
$ws.Range('L6').Value = 20
$ws.Range('L7').Value = 84

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('L3').Value = 49
$ws.Range('L7').Value = 175

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L2').Value = 176
$ws.Range('L3').Value = 174
$ws.Range('L6').Value = 126
$ws.Range('L7').Value = 521

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 97
$ws.Range('L6').Value = 60
$ws.Range('L7').Value = 258

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L3').Value = 67
$ws.Range('L5').Value = 1
$ws.Range('L7').Value = 227

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L3').Value = 329
$ws.Range('L4').Value = 51
$ws.Range('L7').Value = 813

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L6').Value = 86
$ws.Range('L7').Value = 317

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range('L6').Value = 62
$ws.Range('L7').Value = 180

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 305
$ws.Range('L3').Value = 352
$ws.Range('L6').Value = 277
$ws.Range('L7').Value = 1053

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L3').Value = 141
$ws.Range('L7').Value = 351

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 198
$ws.Range('L3').Value = 249
$ws.Range('L4').Value = 44
$ws.Range('L7').Value = 726

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L6').Value = 161
$ws.Range('L7').Value = 583

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 109
$ws.Range('L7').Value = 306

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('L6').Value = 33
$ws.Range('L7').Value = 76

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L2').Value = 158
$ws.Range('L7').Value = 543

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L6').Value = 161
$ws.Range('L7').Value = 331

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 259
$ws.Range('L6').Value = 226
$ws.Range('L7').Value = 866

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L2').Value = 27
$ws.Range('L7').Value = 204

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 153
$ws.Range('L7').Value = 434

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L6').Value = 113
$ws.Range('L7').Value = 245

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('L3').Value = 25
$ws.Range('L7').Value = 72

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L2').Value = 48
$ws.Range('L7').Value = 151

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range('L3').Value = 14
$ws.Range('L7').Value = 42

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L3').Value = 66
$ws.Range('L7').Value = 173

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('L3').Value = 95
$ws.Range('L6').Value = 27
$ws.Range('L7').Value = 216

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('L2').Value = 6
$ws.Range('L6').Value = 26
$ws.Range('L7').Value = 50

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L3').Value = 149
$ws.Range('L6').Value = 90
$ws.Range('L7').Value = 417

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 121
$ws.Range('L6').Value = 106
$ws.Range('L7').Value = 396

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range('L3').Value = 22
$ws.Range('L7').Value = 83

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range('L6').Value = 19
$ws.Range('L7').Value = 77

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L2').Value = 45
$ws.Range('L7').Value = 124

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L2').Value = 45
$ws.Range('L6').Value = 32
$ws.Range('L7').Value = 131

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('L3').Value = 16
$ws.Range('L7').Value = 57

$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('L4').Value = 63
$ws.Range('L7').Value = 115

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range('L6').Value = 4
$ws.Range('L7').Value = 57

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L2').Value = 55
$ws.Range('L7').Value = 195

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('L6').Value = 20
$ws.Range('L7').Value = 60
